$wb = $excel.ActiveWorkbook

# --- Add the two new sheets after the last existing sheet, in order ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSVMZ  = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsSVMZ.Name = "SVM_Z"
$wsSVMZ1 = $wb.Worksheets.Add([Type]::Missing, $wsSVMZ)
$wsSVMZ1.Name = "SVM_Z1"

# Source of the "label" cell format already used throughout the workbook
# (bold font, thin border, centered alignment -> style index 1)
$styleSrc = $wb.Worksheets.Item("XGB_smote3").Range("B1")

# --- Sheet "SVM_Z" : Accuracy / Precision(Avg) / Recall(Avg) ---
$headersZ = @("Accuracy", "Precision(Avg)", "Recall(Avg)")
for ($i = 0; $i -lt $headersZ.Length; $i++) {
    $wsSVMZ.Cells.Item(1, 2 + $i).Value = $headersZ[$i]
}

$dataZ = @(
    @("[0.1, 1]",     0.8202016936184913, 0.8398572549632499, 0.820282038885024),
    @("[0.1, 0.1]",   0.8313459531047886, 0.851020859901092,  0.8314250498125393),
    @("[0.1, 0.01]",  0.8331429606852417, 0.8523256289622969, 0.8332208987570006),
    @("[0.1, 0.001]", 0.7531690989279278, 0.7804823936473813, 0.7532732474442394),
    @("[1, 1]",       0.834118037522918,  0.8445046745608472, 0.8341759138897922),
    @("[1, 0.1]",     0.8337558885243981, 0.8520001665666375, 0.8338319426194308),
    @("[1, 0.01]",    0.8330315080582551, 0.8517541222750274, 0.8331085777033153),
    @("[1, 0.001]",   0.833741960948202,  0.8522439273924404, 0.8338185079466743),
    @("[10, 1]",      0.8262892474949499, 0.8302048856613891, 0.8263254969347413),
    @("[10, 0.1]",    0.8352046261712905, 0.8490241558064276, 0.8352709256667555),
    @("[10, 0.01]",   0.8337698180410438, 0.8522249218111486, 0.8338462777238455),
    @("[10, 0.001]",  0.8329061356168698, 0.8515586436603563, 0.8329830800585452),
    @("[100, 1]",     0.8121639162104985, 0.8133813418809875, 0.8121842487838122),
    @("[100, 0.1]",   0.8362493826096122, 0.8471648237403262, 0.8363086445043949),
    @("[100, 0.01]",  0.8339230495157206, 0.852221501104734,  0.8339991861698973),
    @("[100, 0.001]", 0.8336026541390471, 0.8520668801586948, 0.8336791583789573)
)

for ($r = 0; $r -lt $dataZ.Length; $r++) {
    $row = $r + 2
    $wsSVMZ.Cells.Item($row, 1).Value = $dataZ[$r][0]
    $wsSVMZ.Cells.Item($row, 2).Value = $dataZ[$r][1]
    $wsSVMZ.Cells.Item($row, 3).Value = $dataZ[$r][2]
    $wsSVMZ.Cells.Item($row, 4).Value = $dataZ[$r][3]
}

$styleSrc.Copy()
$wsSVMZ.Range("B1:D1").PasteSpecial(-4122)
$wsSVMZ.Range("A2:A17").PasteSpecial(-4122)

# --- Sheet "SVM_Z1" : Accuracy / Precision(Avg) / Recall(Avg) / std_cross_validation ---
$headersZ1 = @("Accuracy", "Precision(Avg)", "Recall(Avg)", "std_cross_validation")
for ($i = 0; $i -lt $headersZ1.Length; $i++) {
    $wsSVMZ1.Cells.Item(1, 2 + $i).Value = $headersZ1[$i]
}

$dataZ1 = @(
    @("[0.1, 1]",     0.8202016936184913, 0.8398572549632499, 0.820282038885024,  0.002454186493353154),
    @("[0.1, 0.1]",   0.8313459531047886, 0.851020859901092,  0.8314250498125393, 0.002019755326894574),
    @("[0.1, 0.01]",  0.8331429606852417, 0.8523256289622969, 0.8332208987570006, 0.002249333580720828),
    @("[0.1, 0.001]", 0.7531690989279278, 0.7804823936473813, 0.7532732474442394, 0.003187951894907381),
    @("[1, 1]",       0.834118037522918,  0.8445046745608472, 0.8341759138897922, 0.003089174059208984),
    @("[1, 0.1]",     0.8337558885243981, 0.8520001665666375, 0.8338319426194308, 0.002057226086398187),
    @("[1, 0.01]",    0.8330315080582551, 0.8517541222750274, 0.8331085777033153, 0.00202820335174945),
    @("[1, 0.001]",   0.833741960948202,  0.8522439273924404, 0.8338185079466743, 0.002345544505908686),
    @("[10, 1]",      0.8262892474949499, 0.8302048856613891, 0.8263254969347413, 0.002092907043063291),
    @("[10, 0.1]",    0.8352046261712905, 0.8490241558064276, 0.8352709256667555, 0.002501541426657304),
    @("[10, 0.01]",   0.8337698180410438, 0.8522249218111486, 0.8338462777238455, 0.002076162180536619),
    @("[10, 0.001]",  0.8329061356168698, 0.8515586436603563, 0.8329830800585452, 0.002295507980122763),
    @("[100, 1]",     0.8121639162104985, 0.8133813418809875, 0.8121842487838122, 0.0009202182752484279),
    @("[100, 0.1]",   0.8362493826096122, 0.8471648237403262, 0.8363086445043949, 0.003235119973930643),
    @("[100, 0.01]",  0.8339230495157206, 0.852221501104734,  0.8339991861698973, 0.002037621877036707),
    @("[100, 0.001]", 0.8336026541390471, 0.8520668801586948, 0.8336791583789573, 0.00213853529598701)
)

for ($r = 0; $r -lt $dataZ1.Length; $r++) {
    $row = $r + 2
    $wsSVMZ1.Cells.Item($row, 1).Value = $dataZ1[$r][0]
    $wsSVMZ1.Cells.Item($row, 2).Value = $dataZ1[$r][1]
    $wsSVMZ1.Cells.Item($row, 3).Value = $dataZ1[$r][2]
    $wsSVMZ1.Cells.Item($row, 4).Value = $dataZ1[$r][3]
    $wsSVMZ1.Cells.Item($row, 5).Value = $dataZ1[$r][4]
}

$styleSrc.Copy()
$wsSVMZ1.Range("B1:E1").PasteSpecial(-4122)
$wsSVMZ1.Range("A2:A17").PasteSpecial(-4122)

$excel.CutCopyMode = $false
